# Update countries & provincias Spain
# This script refreshes the COVID dashboard "Pais" sheet:
#  - Inserts Mexico (new data) ahead of Peru, shifting Peru/Panama/
#    Republica Dominicana/Islandia down, and drops the old Mexico row,
#    which nets out to a re-ordering of rows 45-50 plus updated stats
#    for Peru and Argentina.
#  - Swaps Camboya ahead of "Consejo Danes para los Refugiados" (rows 117-118)
#  - Swaps Butan ahead of Nepal (rows 196-197)
#  - Refreshes a handful of simple per-country counters (Estados Unidos,
#    Crucero)
#  - Bumps the "last updated" timestamp from 03:50 to 04:20

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "last updated" timestamp -------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 2 de Abril de 2020 a las 04:20"

# --- simple in-place numeric refreshes -----------------------------------
# Estados Unidos (row 4)
$ws.Range("B4").Value = 215086
$ws.Range("C4").Value = 83
$ws.Range("D4").Value = 8878
$ws.Range("E4").Value = 201098
$ws.Range("F4").Value = 5005
$ws.Range("G4").Value = 8
$ws.Range("H4").Value = 5110

# Crucero (row 65) - only Casos activos (D) and Recuperados (E) moved
$ws.Range("D65").Value = 619
$ws.Range("E65").Value = 82

# --- Mexico inserted ahead of Peru, with Peru/Panama/Rep. Dominicana/ -----
# --- Islandia shifted down a row and the old Mexico row removed ----------
$ws.Range("A45").Value = "Mexico"
$ws.Range("B45").Value = 1378
$ws.Range("C45").Value = 163
$ws.Range("D45").Value = 35
$ws.Range("E45").Value = 1306
$ws.Range("F45").Value = 1
$ws.Range("G45").Value = 8
$ws.Range("H45").Value = 37

$ws.Range("A46").Value = "Peru"
$ws.Range("B46").Value = 1323
$ws.Range("C46").Value = 0
$ws.Range("D46").Value = 394
$ws.Range("E46").Value = 882
$ws.Range("F46").Value = 49
$ws.Range("G46").Value = 9
$ws.Range("H46").Value = 47

$ws.Range("A47").Value = "Panama"
$ws.Range("B47").Value = 1317
$ws.Range("C47").Value = 0
$ws.Range("D47").Value = 9
$ws.Range("E47").Value = 1276
$ws.Range("F47").Value = 50
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 32

$ws.Range("A48").Value = "Republica Dominicana"
$ws.Range("B48").Value = 1284
$ws.Range("C48").Value = 0
$ws.Range("D48").Value = 9
$ws.Range("E48").Value = 1218
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 57

$ws.Range("A49").Value = "Islandia"
$ws.Range("B49").Value = 1220
$ws.Range("C49").Value = 0
$ws.Range("D49").Value = 236
$ws.Range("E49").Value = 982
$ws.Range("F49").Value = 12
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 2

$ws.Range("A50").Value = "Argentina"
$ws.Range("B50").Value = 1133
$ws.Range("C50").Value = 0
$ws.Range("D50").Value = 248
$ws.Range("E50").Value = 852
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 1
$ws.Range("H50").Value = 33

# --- Camboya inserted ahead of "Consejo Danes para los Refugiados" -------
$ws.Range("A117").Value = "Camboya"
$ws.Range("B117").Value = 110
$ws.Range("C117").Value = 1
$ws.Range("D117").Value = 34
$ws.Range("E117").Value = 76
$ws.Range("F117").Value = 1
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 0

$ws.Range("A118").Value = "Consejo Danes para los Refugiados"
$ws.Range("B118").Value = 109
$ws.Range("C118").Value = 0
$ws.Range("D118").Value = 3
$ws.Range("E118").Value = 97
$ws.Range("F118").Value = 0
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 9

# --- Butan inserted ahead of Nepal ----------------------------------------
$ws.Range("A196").Value = "Butan"
$ws.Range("B196").Value = 5
$ws.Range("C196").Value = 1
$ws.Range("D196").Value = 1
$ws.Range("E196").Value = 4
$ws.Range("F196").Value = 0
$ws.Range("G196").Value = 0
$ws.Range("H196").Value = 0

$ws.Range("A197").Value = "Nepal"
$ws.Range("B197").Value = 5
$ws.Range("C197").Value = 0
$ws.Range("D197").Value = 1
$ws.Range("E197").Value = 4
$ws.Range("F197").Value = 0
$ws.Range("G197").Value = 0
$ws.Range("H197").Value = 0
